$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D32").Value = "텔레그램봇을 활용한 유저 채팅 데이터 수집 및 활용(feat. telepot, telegram)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/357"

$ws.Range("D37").Value = "[Paper Review] CoST: Contrastive Learning of Disentangled Seasonal-Trend Representations for Time Series Forecasting"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1962&mod=document&pageid=1"

$ws.Range("D46").Value = "호흡곤란(숨참) 감별진단"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/442"
